$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.9591836734693877
$ws.Range("D2").Value = 0.9791666666666666
$ws.Range("E2").Value = 49

$ws.Range("B3").Value = 0.9565217391304348
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.9777777777777777
$ws.Range("E3").Value = 44

$ws.Range("B4").Value = 0.978494623655914
$ws.Range("C4").Value = 0.978494623655914
$ws.Range("D4").Value = 0.978494623655914
$ws.Range("E4").Value = 0.978494623655914

$ws.Range("B5").Value = 0.9782608695652174
$ws.Range("C5").Value = 0.9795918367346939
$ws.Range("D5").Value = 0.9784722222222222
$ws.Range("E5").Value = 93

$ws.Range("B6").Value = 0.9794296400187003
$ws.Range("C6").Value = 0.978494623655914
$ws.Range("D6").Value = 0.9785095579450418
$ws.Range("E6").Value = 93

$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 0.9795918367346939
$ws.Range("D7").Value = 0.9896907216494846
$ws.Range("E7").Value = 49

$ws.Range("B8").Value = 0.9777777777777777
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0.9887640449438202
$ws.Range("E8").Value = 44

$ws.Range("B9").Value = 0.989247311827957
$ws.Range("C9").Value = 0.989247311827957
$ws.Range("D9").Value = 0.989247311827957
$ws.Range("E9").Value = 0.989247311827957

$ws.Range("B10").Value = 0.9888888888888889
$ws.Range("C10").Value = 0.9897959183673469
$ws.Range("D10").Value = 0.9892273832966524
$ws.Range("E10").Value = 93

$ws.Range("B11").Value = 0.9894862604540022
$ws.Range("C11").Value = 0.989247311827957
$ws.Range("D11").Value = 0.9892522939607831
$ws.Range("E11").Value = 93

$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 0.9795918367346939
$ws.Range("D12").Value = 0.9896907216494846
$ws.Range("E12").Value = 49

$ws.Range("B13").Value = 0.9777777777777777
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 0.9887640449438202
$ws.Range("E13").Value = 44

$ws.Range("B14").Value = 0.989247311827957
$ws.Range("C14").Value = 0.989247311827957
$ws.Range("D14").Value = 0.989247311827957
$ws.Range("E14").Value = 0.989247311827957

$ws.Range("B15").Value = 0.9888888888888889
$ws.Range("C15").Value = 0.9897959183673469
$ws.Range("D15").Value = 0.9892273832966524
$ws.Range("E15").Value = 93

$ws.Range("B16").Value = 0.9894862604540022
$ws.Range("C16").Value = 0.989247311827957
$ws.Range("D16").Value = 0.9892522939607831
$ws.Range("E16").Value = 93

$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 0.9387755102040817
$ws.Range("D17").Value = 0.968421052631579
$ws.Range("E17").Value = 49

$ws.Range("B18").Value = 0.9361702127659575
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 0.967032967032967
$ws.Range("E18").Value = 44

$ws.Range("B19").Value = 0.967741935483871
$ws.Range("C19").Value = 0.967741935483871
$ws.Range("D19").Value = 0.967741935483871
$ws.Range("E19").Value = 0.967741935483871

$ws.Range("B20").Value = 0.9680851063829787
$ws.Range("C20").Value = 0.9693877551020409
$ws.Range("D20").Value = 0.967727009832273
$ws.Range("E20").Value = 93

$ws.Range("B21").Value = 0.9698009608785175
$ws.Range("C21").Value = 0.967741935483871
$ws.Range("D21").Value = 0.9677643239612679
$ws.Range("E21").Value = 93

$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 0.9387755102040817
$ws.Range("D22").Value = 0.968421052631579
$ws.Range("E22").Value = 49

$ws.Range("B23").Value = 0.9361702127659575
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 0.967032967032967
$ws.Range("E23").Value = 44

$ws.Range("B24").Value = 0.967741935483871
$ws.Range("C24").Value = 0.967741935483871
$ws.Range("D24").Value = 0.967741935483871
$ws.Range("E24").Value = 0.967741935483871

$ws.Range("B25").Value = 0.9680851063829787
$ws.Range("C25").Value = 0.9693877551020409
$ws.Range("D25").Value = 0.967727009832273
$ws.Range("E25").Value = 93

$ws.Range("B26").Value = 0.9698009608785175
$ws.Range("C26").Value = 0.967741935483871
$ws.Range("D26").Value = 0.9677643239612679
$ws.Range("E26").Value = 93
